$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 355.56668
$ws.Range("J17").Value = 387.57693
$ws.Range("L17").Value = 1162.73079
$ws.Range("N17").Value = -1498.73079
$ws.Range("H39").Value = 2008.1666
$ws.Range("I39").Value = 1509.125
$ws.Range("J39").Value = 3006.25
$ws.Range("K39").Value = 4527.375
$ws.Range("L39").Value = 9018.75
$ws.Range("M39").Value = -4231.375
$ws.Range("N39").Value = -9610.75
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("L52").ClearContents()
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = 0
$ws.Range("H62").Value = 66697300
$ws.Range("I62").Value = 111112344
$ws.Range("K62").Value = 111112344
$ws.Range("M62").Value = -111111720
$ws.Range("H65").Value = 66697300
$ws.Range("I65").Value = 111112344
$ws.Range("K65").Value = 555561720
$ws.Range("M65").Value = -555558600
$ws.Range("H69").Value = 4994
$ws.Range("I69").Value = 4994
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 14982
$ws.Range("L69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -14108
$ws.Range("H72").Value = 4994
$ws.Range("I72").Value = 4994
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 44946
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -40578
$ws.Range("N72").Value = 0
$ws.Range("H103").Value = 581.9048
$ws.Range("J103").Value = 582.3889
$ws.Range("L103").Value = 1747.1667
$ws.Range("N103").Value = -2919.1667
$ws.Range("H112").Value = 5663.795
$ws.Range("J112").Value = 6059.6113
$ws.Range("L112").Value = 18178.8339
$ws.Range("N112").Value = -20394.8339
$ws.Range("H124").Value = 101300
$ws.Range("J124").Value = 101950
$ws.Range("L124").Value = 101950
$ws.Range("N124").Value = -111770
$ws.Range("I125").Value = 50000776
$ws.Range("K125").Value = 450006984
$ws.Range("M125").Value = -450004524
$ws.Range("H131").Value = 1917.375
$ws.Range("J131").Value = 2405.5
$ws.Range("L131").Value = 7216.5
$ws.Range("N131").Value = -17296.5
$ws.Range("H135").Value = 589110.1
$ws.Range("J135").Value = 2888.3333
$ws.Range("L135").Value = 25994.9997
$ws.Range("N135").Value = -31064.9997
$ws.Range("H137").Value = 8073.5713
$ws.Range("I137").Value = 7804.3
$ws.Range("J137").Value = 8746.75
$ws.Range("K137").Value = 23412.9
$ws.Range("L137").Value = 26240.25
$ws.Range("M137").Value = -20862.9
$ws.Range("N137").Value = -31340.25
$ws.Range("H138").Value = 5251.3706
$ws.Range("I138").Value = 2144.353
$ws.Range("J138").Value = 6678.919
$ws.Range("K138").Value = 6433.059
$ws.Range("L138").Value = 20036.757
$ws.Range("M138").Value = -1293.059
$ws.Range("N138").Value = -30316.757

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5212.1924
$ws.Range("I61").Value = 3640.111
$ws.Range("J61").Value = 8749.375
$ws.Range("K61").Value = 3640.111
$ws.Range("L61").Value = 8749.375
$ws.Range("M61").Value = -3428.111
$ws.Range("N61").Value = -9173.375
$ws.Range("H74").Value = 184779.11
$ws.Range("J74").Value = 15338
$ws.Range("L74").Value = 15338
$ws.Range("N74").Value = -17086
$ws.Range("H77").Value = 184779.11
$ws.Range("J77").Value = 15338
$ws.Range("L77").Value = 76690
$ws.Range("N77").Value = -85426
$ws.Range("H97").Value = 2878100.5
$ws.Range("I97").Value = 338.2857
$ws.Range("K97").Value = 338.2857
$ws.Range("M97").Value = 157.7143
$ws.Range("H132").Value = 7424.3477
$ws.Range("I132").Value = 7647.5
$ws.Range("K132").Value = 22942.5
$ws.Range("M132").Value = -20412.5
$ws.Range("H136").Value = 5212.1924
$ws.Range("I136").Value = 3640.111
$ws.Range("J136").Value = 8749.375
$ws.Range("K136").Value = 10920.333
$ws.Range("L136").Value = 26248.125
$ws.Range("M136").Value = -8370.332999999999
$ws.Range("N136").Value = -31348.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("H20").Value = 9260908
$ws.Range("I20").Value = 15153206
$ws.Range("J20").Value = 1581.8572
$ws.Range("K20").Value = 15153206
$ws.Range("L20").Value = 1581.8572
$ws.Range("M20").Value = -15152959
$ws.Range("N20").Value = -2075.8572
$ws.Range("H99").Value = 4331543.5
$ws.Range("I99").Value = 2569.842
$ws.Range("K99").Value = 2569.842
$ws.Range("M99").Value = -1071.842
$ws.Range("H105").Value = 3785.9048
$ws.Range("I105").Value = 2428.5
$ws.Range("K105").Value = 2428.5
$ws.Range("M105").Value = -681.5
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").ClearContents()
$ws.Range("N133").Value = 0
$ws.Range("H134").Value = 4347.5454
$ws.Range("I134").Value = 3198.3225
$ws.Range("J134").Value = 7088
$ws.Range("K134").Value = 9594.967500000001
$ws.Range("L134").Value = 21264
$ws.Range("M134").Value = -7059.967500000001
$ws.Range("N134").Value = -26334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15906.546
$ws.Range("I31").Value = 7748.25
$ws.Range("J31").Value = 20568.428
$ws.Range("K31").Value = 7748.25
$ws.Range("L31").Value = 20568.428
$ws.Range("M31").Value = -7453.25
$ws.Range("N31").Value = -21158.428
$ws.Range("H34").Value = 15906.546
$ws.Range("I34").Value = 7748.25
$ws.Range("J34").Value = 20568.428
$ws.Range("K34").Value = 7748.25
$ws.Range("L34").Value = 20568.428
$ws.Range("M34").Value = -7546.25
$ws.Range("N34").Value = -20972.428
$ws.Range("H107").Value = 2097.8572
$ws.Range("I107").Value = 984.375
$ws.Range("J107").Value = 3582.5
$ws.Range("K107").Value = 984.375
$ws.Range("L107").Value = 3582.5
$ws.Range("M107").Value = 935.625
$ws.Range("N107").Value = -7422.5
$ws.Range("H108").Value = 60652.332
$ws.Range("J108").Value = 60652.332
$ws.Range("L108").Value = 60652.332
$ws.Range("N108").Value = -68332.33199999999
$ws.Range("H134").Value = 3224.1667
$ws.Range("I134").Value = 2374.7083
$ws.Range("K134").Value = 7124.124899999999
$ws.Range("M134").Value = -4589.124899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 223420.22
$ws.Range("I2").Value = 130
$ws.Range("K2").Value = 780
$ws.Range("M2").Value = -667
$ws.Range("I12").Value = 649.7143
$ws.Range("J12").Value = 3846787.8
$ws.Range("K12").Value = 1949.1429
$ws.Range("L12").Value = 11540363.4
$ws.Range("M12").Value = -1776.1429
$ws.Range("N12").Value = -11540709.4
$ws.Range("H26").Value = 266.33334
$ws.Range("I26").Value = 98
$ws.Range("K26").Value = 294
$ws.Range("M26").Value = -6
$ws.Range("H109").Value = 67904110
$ws.Range("J109").Value = 22225802
$ws.Range("L109").Value = 66677406
$ws.Range("N109").Value = -66679486
$ws.Range("H115").Value = 918.7143
$ws.Range("I115").Value = 646.2
$ws.Range("J115").Value = 1600
$ws.Range("K115").Value = 1938.6
$ws.Range("L115").Value = 4800
$ws.Range("M115").Value = -763.6000000000001
$ws.Range("N115").Value = -7150
$ws.Range("H119").Value = 4684.6665
$ws.Range("I119").Value = 4684.6665
$ws.Range("K119").Value = 14053.9995
$ws.Range("M119").Value = -9215.999500000002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 3431
$ws.Range("J13").Value = 4999
$ws.Range("L13").Value = 4999
$ws.Range("N13").Value = -5277
$ws.Range("H57").Value = 50500
$ws.Range("H80").Value = 255198.75
$ws.Range("I80").Value = 6900
$ws.Range("K80").Value = 6900
$ws.Range("M80").Value = -5902
$ws.Range("H83").Value = 255198.75
$ws.Range("I83").Value = 6900
$ws.Range("K83").Value = 34500
$ws.Range("M83").Value = -29508
$ws.Range("H97").Value = 1377.303
$ws.Range("I97").Value = 1194.3077
$ws.Range("K97").Value = 1194.3077
$ws.Range("M97").Value = -698.3077000000001
$ws.Range("H122").Value = 24446.695
$ws.Range("I122").Value = 30695.086
$ws.Range("J122").Value = 4565.4546
$ws.Range("K122").Value = 92085.258
$ws.Range("L122").Value = 13696.3638
$ws.Range("M122").Value = -89635.258
$ws.Range("N122").Value = -18596.3638

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2464166.8
$ws.Range("J46").Value = 1257.8182
$ws.Range("L46").Value = 1257.8182
$ws.Range("N46").Value = -1633.8182
$ws.Range("H109").Value = 55000
$ws.Range("J109").Value = 55000
$ws.Range("L109").Value = 55000
$ws.Range("N109").Value = -57774
$ws.Range("H136").Value = 5287.5
$ws.Range("I136").Value = 2875
$ws.Range("K136").Value = 8625
$ws.Range("M136").Value = -6075

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2455.4443
$ws.Range("I96").Value = 2380
$ws.Range("J96").Value = 2549.75
$ws.Range("K96").Value = 2380
$ws.Range("L96").Value = 2549.75
$ws.Range("M96").Value = -1007
$ws.Range("N96").Value = -5295.75
$ws.Range("H122").Value = 19390830
$ws.Range("I122").Value = 24005646
$ws.Range("K122").Value = 72016938
$ws.Range("M122").Value = -72014488
$ws.Range("H132").Value = 35759452
$ws.Range("I132").Value = 55569316
$ws.Range("J132").Value = 101700
$ws.Range("K132").Value = 166707948
$ws.Range("L132").Value = 305100
$ws.Range("M132").Value = -166705418
$ws.Range("N132").Value = -310160
